# Scheduled runner update: refresh cached market-board price/profit figures
# across the per-job Leve profit sheets (currentAveragePrice*, LevePrice*,
# LeveProfit* columns). Mirrors the data the external pricing job writes
# back into the workbook on each run.
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(74, 8).Value = 4258.8335
$ws.Cells.Item(74, 9).Value = 4258.8335
$ws.Cells.Item(74, 11).Value = 4258.8335
$ws.Cells.Item(74, 13).Value = -3322.8335
$ws.Cells.Item(77, 8).Value = 4258.8335
$ws.Cells.Item(77, 9).Value = 4258.8335
$ws.Cells.Item(77, 11).Value = 21294.1675
$ws.Cells.Item(77, 13).Value = -16614.1675
$ws.Cells.Item(138, 8).Value = 18982.553
$ws.Cells.Item(138, 9).Value = 823.5946
$ws.Cells.Item(138, 10).Value = 54344.74
$ws.Cells.Item(138, 11).Value = 2470.7838
$ws.Cells.Item(138, 12).Value = 163034.22
$ws.Cells.Item(138, 13).Value = 2669.2162
$ws.Cells.Item(138, 14).Value = -173314.22

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(3, 8).Value = 36168
$ws.Cells.Item(3, 10).Value = 53754.5
$ws.Cells.Item(3, 12).Value = 53754.5
$ws.Cells.Item(3, 14).Value = -53984.5
$ws.Cells.Item(110, 8).Value = 5466.6665
$ws.Cells.Item(110, 9).Value = 1933.3334
$ws.Cells.Item(110, 10).Value = 9000
$ws.Cells.Item(110, 11).Value = 1933.3334
$ws.Cells.Item(110, 12).Value = 9000
$ws.Cells.Item(110, 13).Value = 111.6666
$ws.Cells.Item(110, 14).Value = -13090
$ws.Cells.Item(132, 8).Value = 9018.634
$ws.Cells.Item(132, 9).Value = 8277.9375
$ws.Cells.Item(132, 10).Value = 9865.143
$ws.Cells.Item(132, 11).Value = 24833.8125
$ws.Cells.Item(132, 12).Value = 29595.429
$ws.Cells.Item(132, 13).Value = -22303.8125
$ws.Cells.Item(132, 14).Value = -34655.429

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(64, 8).Value = 499
$ws.Cells.Item(64, 9).Value = 675.3333
$ws.Cells.Item(64, 10).Value = 417.6154
$ws.Cells.Item(64, 11).Value = 675.3333
$ws.Cells.Item(64, 12).Value = 417.6154
$ws.Cells.Item(64, 13).Value = -450.3333
$ws.Cells.Item(64, 14).Value = -867.6154
$ws.Cells.Item(67, 8).Value = 499
$ws.Cells.Item(67, 9).Value = 675.3333
$ws.Cells.Item(67, 10).Value = 417.6154
$ws.Cells.Item(67, 11).Value = 675.3333
$ws.Cells.Item(67, 12).Value = 417.6154
$ws.Cells.Item(67, 13).Value = 104.6667
$ws.Cells.Item(67, 14).Value = -1977.6154

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2124.64
$ws.Cells.Item(31, 9).Value = 1609.5151
$ws.Cells.Item(31, 10).Value = 3124.5881
$ws.Cells.Item(31, 11).Value = 1609.5151
$ws.Cells.Item(31, 12).Value = 3124.5881
$ws.Cells.Item(31, 13).Value = -1314.5151
$ws.Cells.Item(31, 14).Value = -3714.5881
$ws.Cells.Item(34, 8).Value = 2124.64
$ws.Cells.Item(34, 9).Value = 1609.5151
$ws.Cells.Item(34, 10).Value = 3124.5881
$ws.Cells.Item(34, 11).Value = 1609.5151
$ws.Cells.Item(34, 12).Value = 3124.5881
$ws.Cells.Item(34, 13).Value = -1407.5151
$ws.Cells.Item(34, 14).Value = -3528.5881
$ws.Cells.Item(132, 8).Value = 3127405.5
$ws.Cells.Item(132, 9).Value = 5557748.5
$ws.Cells.Item(132, 10).Value = 2678.9285
$ws.Cells.Item(132, 11).Value = 16673245.5
$ws.Cells.Item(132, 12).Value = 8036.7855
$ws.Cells.Item(132, 13).Value = -16670715.5
$ws.Cells.Item(132, 14).Value = -13096.7855
$ws.Cells.Item(134, 8).Value = 3198.3914
$ws.Cells.Item(134, 9).Value = 1821.2222
$ws.Cells.Item(134, 10).Value = 4083.7144
$ws.Cells.Item(134, 11).Value = 5463.6666
$ws.Cells.Item(134, 12).Value = 12251.1432
$ws.Cells.Item(134, 13).Value = -2928.6666
$ws.Cells.Item(134, 14).Value = -17321.1432

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 943.6316
$ws.Cells.Item(5, 9).Value = 550.2727
$ws.Cells.Item(5, 10).Value = 1484.5
$ws.Cells.Item(5, 11).Value = 1650.8181
$ws.Cells.Item(5, 12).Value = 4453.5
$ws.Cells.Item(5, 13).Value = -1538.8181
$ws.Cells.Item(5, 14).Value = -4677.5
$ws.Cells.Item(131, 8).Value = 532307.9
$ws.Cells.Item(131, 9).Value = 486.8125
$ws.Cells.Item(131, 10).Value = 1064129
$ws.Cells.Item(131, 11).Value = 1460.4375
$ws.Cells.Item(131, 12).Value = 3192387
$ws.Cells.Item(131, 13).Value = 3579.5625
$ws.Cells.Item(131, 14).Value = -3202467
$ws.Cells.Item(135, 8).Value = 943.6316
$ws.Cells.Item(135, 9).Value = 550.2727
$ws.Cells.Item(135, 10).Value = 1484.5
$ws.Cells.Item(135, 11).Value = 4952.454299999999
$ws.Cells.Item(135, 12).Value = 13360.5
$ws.Cells.Item(135, 13).Value = -2417.454299999999
$ws.Cells.Item(135, 14).Value = -18430.5

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 9878.615
$ws.Cells.Item(132, 9).Value = 21284.8
$ws.Cells.Item(132, 10).Value = 2749.75
$ws.Cells.Item(132, 11).Value = 63854.39999999999
$ws.Cells.Item(132, 12).Value = 8249.25
$ws.Cells.Item(132, 13).Value = -61324.39999999999
$ws.Cells.Item(132, 14).Value = -13309.25

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(32, 8).Value = 11500
$ws.Cells.Item(32, 9).Value = 11500
$ws.Cells.Item(32, 10).Value = 0
$ws.Cells.Item(32, 11).Value = 11500
$ws.Cells.Item(32, 12).Value = $null
$ws.Cells.Item(32, 13).Value = -11183
$ws.Cells.Item(32, 14).Value = 0
$ws.Cells.Item(33, 8).Value = 56678
$ws.Cells.Item(33, 10).Value = 70017
$ws.Cells.Item(33, 12).Value = 70017
$ws.Cells.Item(33, 14).Value = -70597
$ws.Cells.Item(34, 8).Value = 0
$ws.Cells.Item(34, 9).Value = 0
$ws.Cells.Item(34, 11).Value = 0
$ws.Cells.Item(34, 13).Value = $null
$ws.Cells.Item(45, 8).Value = 10000
$ws.Cells.Item(45, 10).Value = 10000
$ws.Cells.Item(45, 12).Value = 10000
$ws.Cells.Item(45, 14).Value = -10814
$ws.Cells.Item(46, 8).Value = 834294.9399999999
$ws.Cells.Item(46, 9).Value = 750
$ws.Cells.Item(46, 10).Value = 3334929.8
$ws.Cells.Item(46, 11).Value = 750
$ws.Cells.Item(46, 12).Value = 3334929.8
$ws.Cells.Item(46, 13).Value = -562
$ws.Cells.Item(46, 14).Value = -3335305.8
$ws.Cells.Item(56, 8).Value = 8000
$ws.Cells.Item(56, 9).Value = 8000
$ws.Cells.Item(56, 10).Value = 0
$ws.Cells.Item(56, 11).Value = 8000
$ws.Cells.Item(56, 12).Value = 0
$ws.Cells.Item(56, 13).Value = $null
$ws.Cells.Item(56, 14).Value = -7309
$ws.Cells.Item(58, 8).Value = 4980
$ws.Cells.Item(58, 9).Value = 4980
$ws.Cells.Item(58, 10).Value = 0
$ws.Cells.Item(58, 11).Value = 4980
$ws.Cells.Item(58, 12).Value = 0
$ws.Cells.Item(58, 13).Value = $null
$ws.Cells.Item(58, 14).Value = -4720
$ws.Cells.Item(93, 8).Value = 1807.9412
$ws.Cells.Item(93, 9).Value = 1627.5
$ws.Cells.Item(93, 10).Value = 2650
$ws.Cells.Item(93, 11).Value = 1627.5
$ws.Cells.Item(93, 12).Value = 2650
$ws.Cells.Item(93, 13).Value = -379.5
$ws.Cells.Item(93, 14).Value = -5146
$ws.Cells.Item(132, 8).Value = 43481772
$ws.Cells.Item(132, 9).Value = 58825280
$ws.Cells.Item(132, 10).Value = 8499.666999999999
$ws.Cells.Item(132, 11).Value = 176475840
$ws.Cells.Item(132, 12).Value = 25499.001
$ws.Cells.Item(132, 13).Value = -176473310
$ws.Cells.Item(132, 14).Value = -30559.001

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 115573736
$ws.Cells.Item(2, 10).Value = 200031520
$ws.Cells.Item(2, 12).Value = 200031520
$ws.Cells.Item(2, 14).Value = -200031744
$ws.Cells.Item(81, 8).Value = 1445.5
$ws.Cells.Item(81, 10).Value = 1660.6666
$ws.Cells.Item(81, 12).Value = 3321.3332
$ws.Cells.Item(81, 14).Value = -5443.3332
$ws.Cells.Item(84, 8).Value = 1445.5
$ws.Cells.Item(84, 10).Value = 1660.6666
$ws.Cells.Item(84, 12).Value = 16606.666
$ws.Cells.Item(84, 14).Value = -27214.666
$ws.Cells.Item(132, 8).Value = 3082.0715
$ws.Cells.Item(132, 9).Value = 2441.3333
$ws.Cells.Item(132, 10).Value = 3562.625
$ws.Cells.Item(132, 11).Value = 7323.999899999999
$ws.Cells.Item(132, 12).Value = 10687.875
$ws.Cells.Item(132, 13).Value = -4793.999899999999
$ws.Cells.Item(132, 14).Value = -15747.875
$ws.Cells.Item(136, 8).Value = 6251002
$ws.Cells.Item(136, 9).Value = 7353831.5
$ws.Cells.Item(136, 10).Value = 1633.3334
$ws.Cells.Item(136, 11).Value = 22061494.5
$ws.Cells.Item(136, 12).Value = 4900.0002
$ws.Cells.Item(136, 13).Value = -22058944.5
$ws.Cells.Item(136, 14).Value = -10000.0002
